{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values && table.values.length ? table.values[0].length : 4;\n\n// Center every cell in the table vertically.\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.verticalAlignment = Word.VerticalAlignment.center;\n  }\n}\n\n// The very first cell (row 1, column 1) also shrinks its \"${barcode}\" text\n// (and its paragraph mark) from 10.5pt to 9.5pt (sz 21 -> 19, szCs -> 20).\nconst firstCell = table.getCell(0, 0);\nfirstCell.body.font.size = 9.5;\nfirstCell.body.font.sizeBidirectional = 10;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n# Center every cell in the table vertically.\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.VerticalAlignment = 1\n    }\n}\n\n# The very first cell (row 1, column 1) also shrinks its \"${barcode}\" text\n# (and its paragraph mark) from 10.5pt to 9.5pt (sz 21 -> 19, szCs -> 20).\n$firstCell = $table.Cell(1, 1)\n$firstCell.Range.Font.Size = 9.5\n$firstCell.Range.Font.SizeBi = 10\n"}
